# Natmi LR-pairs table update ("Natmi following Dr Hou advice"):
# recompute rows for the ECs/FAPs/sCs sending/target cluster combinations
# (Ligand=Reln, Receptor=Vldlr), expanding the table from 4 to 9 data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @()
$rows += ,@("ECs","Reln","Vldlr","ECs",2,0.6666666666666666,0.2638773333333334,0.791632,0.0472992197617186,0.04729921976171859,2,0.6666666666666666,0.3374003333333333,1.012201,0.01738364872808817,0.01738364872808818,0.08903230022577778,0.8012907020319999,0.0008222330214503627,0.0008222330214503627)
$rows += ,@("ECs","Reln","Vldlr","FAPs",2,0.6666666666666666,0.2638773333333334,0.791632,0.0472992197617186,0.04729921976171859,3,1,16.71131166666667,50.133935,0.8610055862391021,0.8610055862391023,4.409736359102222,39.68762723192,0.04072489243959065,0.04072489243959065)
$rows += ,@("ECs","Reln","Vldlr","sCs",2,0.6666666666666666,0.2638773333333334,0.791632,0.0472992197617186,0.04729921976171859,3,1,2.360351,7.081053000000001,0.1216107650328097,0.1216107650328097,0.6228431276106667,5.605588148496,0.005752094300677588,0.005752094300677589)
$rows += ,@("FAPs","Reln","Vldlr","ECs",3,1,1.127914,3.383742,0.2021751981665183,0.2021751981665183,2,0.6666666666666666,0.3374003333333333,1.012201,0.01738364872808817,0.01738364872808818,0.3805585595713333,3.425027036142,0.00351454262645837,0.003514542626458371)
$rows += ,@("FAPs","Reln","Vldlr","FAPs",3,1,1.127914,3.383742,0.2021751981665183,0.2021751981665183,3,1,16.71131166666667,50.133935,0.8610055862391021,0.8610055862391023,18.84892238719667,169.64030148477,0.1740739750203697,0.1740739750203698)
$rows += ,@("FAPs","Reln","Vldlr","sCs",3,1,1.127914,3.383742,0.2021751981665183,0.2021751981665183,3,1,2.360351,7.081053000000001,0.1216107650328097,0.1216107650328097,2.662272937814,23.960456440326,0.02458668051969018,0.02458668051969019)
$rows += ,@("sCs","Reln","Vldlr","ECs",3,1,4.187102666666667,12.561308,0.7505255820717631,0.750525582071763,2,0.6666666666666666,0.3374003333333333,1.012201,0.01738364872808817,0.01738364872808818,1.412729835434222,12.714568518908,0.01304687308017944,0.01304687308017944)
$rows += ,@("sCs","Reln","Vldlr","FAPs",3,1,4.187102666666667,12.561308,0.7505255820717631,0.750525582071763,3,1,16.71131166666667,50.133935,0.8610055862391021,0.8610055862391023,69.97197764299779,629.7477987869801,0.6462067187791417,0.6462067187791418)
$rows += ,@("sCs","Reln","Vldlr","sCs",3,1,4.187102666666667,12.561308,0.7505255820717631,0.750525582071763,3,1,2.360351,7.081053000000001,0.1216107650328097,0.1216107650328097,9.883031966369336,88.94728769732401,0.09127199021244188,0.0912719902124419)

for ($i = 0; $i -lt $rows.Length; $i++) {
  $r = $i + 2
  $row = $rows[$i]
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  for ($j = 0; $j -lt 16; $j++) {
    $ws.Cells.Item($r, 5 + $j).Value = $row[4 + $j]
  }
}
